$d = $word.ActiveDocument

# Locate the date text "11/05/2022" in the document.
$found = $d.Content
$ok = $found.Find.Execute("11/05/2022", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($ok) {
    # Narrow the range down to just the day portion ("11") that changes to "16",
    # leaving the rest of the date ("/05/2022") untouched.
    $dayRange = $d.Range($found.Start, $found.Start + 2)
    $dayRange.Text = "16"

    # Re-apply (identical) direct character formatting to the replaced portion so
    # that it is stored as its own run, distinct from the remainder of the date.
    $dayRange2 = $d.Range($found.Start, $found.Start + 2)
    $dayRange2.Font.Name = "Segoe UI"
    $dayRange2.Font.NameBi = "Segoe UI"
}
